$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1404
$ws1.Range("F4").Value = 87
$ws1.Range("F7").Value = 527
$ws1.Range("F8").Value = 859
$ws1.Range("F9").Value = 585
$ws1.Range("F10").Value = 778
$ws1.Range("F11").Value = 350
$ws1.Range("F12").Value = 530
$ws1.Range("F13").Value = 102
$ws1.Range("F14").Value = 1103
$ws1.Range("F15").Value = 541
$ws1.Range("F16").Value = 316
$ws1.Range("F18").Value = 129
$ws1.Range("F19").Value = 279
$ws1.Range("F20").Value = 41
$ws1.Range("F22").Value = 514
$ws1.Range("F23").Value = 501
$ws1.Range("F24").Value = 6
$ws1.Range("F25").Value = 423

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 396
$ws2.Range("F3").Value = 63
$ws2.Range("F5").Value = 302
$ws2.Range("F6").Value = 94
$ws2.Range("F9").Value = 168
$ws2.Range("F10").Value = 177

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1404
$ws4.Range("F5").Value = 87
$ws4.Range("F6").Value = 396
$ws4.Range("F8").Value = 63
$ws4.Range("F11").Value = 302
$ws4.Range("F12").Value = 527
$ws4.Range("F13").Value = 859
$ws4.Range("F14").Value = 585
$ws4.Range("F15").Value = 778
$ws4.Range("F16").Value = 350
$ws4.Range("F17").Value = 530
$ws4.Range("F18").Value = 102
$ws4.Range("F19").Value = 1103
$ws4.Range("F20").Value = 541
$ws4.Range("F21").Value = 94
$ws4.Range("F23").Value = 316
$ws4.Range("F26").Value = 129
$ws4.Range("F27").Value = 168
$ws4.Range("F28").Value = 279
$ws4.Range("F29").Value = 41
$ws4.Range("F31").Value = 177
$ws4.Range("F33").Value = 514
$ws4.Range("F36").Value = 501
$ws4.Range("F37").Value = 6
$ws4.Range("F38").Value = 423
